$wb = $excel.ActiveWorkbook

# --- Update "Conversión del día" message on Hoja1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.06 = 24223.38 pesos`n✅ 24223.38 pesos = 6.04 = 947.92 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update tasas (rates) values on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 165.01
$ws2.Range("O10").Value = 3997.1
$ws2.Range("N12").Value = 4012
$ws2.Range("O12").Value = 157
